# "Carico Modello Fisico + alcune modifiche dd"
#
# Data-dictionary maintenance pass:
#   - The ID_Corso foreign-key "VARCHAR(10)" data type was wrong; it is
#     corrected to "CHAR(10)" everywhere it appears (PRESENZA!D3, FIRMA!D7,
#     LEZIONE!D7).
#   - A few sheet selections were left in different cells while reviewing,
#     and the LEZIONE sheet ends up as the active tab.

$wb = $excel.ActiveWorkbook

$wsRegistro = $wb.Worksheets.Item("REGISTRO")
$wsPresenza = $wb.Worksheets.Item("PRESENZA")
$wsFirma    = $wb.Worksheets.Item("FIRMA")
$wsStudente = $wb.Worksheets.Item("STUDENTE")
$wsDocente  = $wb.Worksheets.Item("DOCENTE")
$wsLezione  = $wb.Worksheets.Item("LEZIONE")

# --- Data type correction: VARCHAR(10) -> CHAR(10) ------------------------
$wsPresenza.Range("D3").Value = "CHAR(10)"
$wsFirma.Range("D7").Value = "CHAR(10)"
$wsLezione.Range("D7").Value = "CHAR(10)"

# --- Per-sheet selection bookkeeping --------------------------------------
$wsRegistro.Activate()
$wsRegistro.Range("D10").Select()

$wsFirma.Activate()
$wsFirma.Range("D7").Select()

$wsLezione.Activate()
$wsLezione.Range("D7").Select()

# LEZIONE is the sheet left active/visible when the workbook is saved.
$wsLezione.Activate()
